# Auto-applies the numeric cell updates described by the commit diff
# (Titan_Profits leve-profit recompute across ALC/ARM/BSM/CRP/CUL/GSM/LTW).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1339.4546
$ws.Range("J129").Value = 2480
$ws.Range("L129").Value = 7440
$ws.Range("N129").Value = -17440

$ws.Range("H137").Value = 20409914
$ws.Range("I137").Value = 34483776
$ws.Range("J137").Value = 2815.15
$ws.Range("K137").Value = 103451328
$ws.Range("L137").Value = 8445.450000000001
$ws.Range("M137").Value = -103448778
$ws.Range("N137").Value = -13545.45

$ws.Range("H138").Value = 2095.32
$ws.Range("I138").Value = 1025
$ws.Range("J138").Value = 2779.623
$ws.Range("K138").Value = 3075
$ws.Range("L138").Value = 8338.869000000001
$ws.Range("M138").Value = 2065
$ws.Range("N138").Value = -18618.869

$ws.Range("H141").Value = 5855.75
$ws.Range("I141").Value = 2953.182
$ws.Range("J141").Value = 9403.333000000001
$ws.Range("K141").Value = 8859.545999999998
$ws.Range("L141").Value = 28209.999
$ws.Range("M141").Value = -3679.545999999998
$ws.Range("N141").Value = -38569.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 550
$ws.Range("I3").Value = 550
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 550
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -435
$ws.Range("N3").Value = ""

$ws.Range("H32").Value = 16699.135
$ws.Range("I32").Value = 4660.091
$ws.Range("J32").Value = 134414.22
$ws.Range("K32").Value = 4660.091
$ws.Range("L32").Value = 134414.22
$ws.Range("M32").Value = -4373.091
$ws.Range("N32").Value = -134988.22

$ws.Range("H45").Value = 1500
$ws.Range("I45").Value = 1750
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 1750
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -1373
$ws.Range("N45").Value = -1754

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 59000
$ws.Range("J59").Value = 59000
$ws.Range("L59").Value = 59000
$ws.Range("N59").Value = -60694

$ws.Range("H75").Value = 800236
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 800236
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 800236
$ws.Range("M75").Value = ""
$ws.Range("N75").Value = -802108

$ws.Range("H78").Value = 800236
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 800236
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 2400708
$ws.Range("M78").Value = ""
$ws.Range("N78").Value = -2410068

$ws.Range("H135").Value = 49070
$ws.Range("J135").Value = 49070
$ws.Range("L135").Value = 49070
$ws.Range("N135").Value = -59210

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3838.7646
$ws.Range("I31").Value = 1792.8182
$ws.Range("J31").Value = 5767.8
$ws.Range("K31").Value = 1792.8182
$ws.Range("L31").Value = 5767.8
$ws.Range("M31").Value = -1497.8182
$ws.Range("N31").Value = -6357.8

$ws.Range("H34").Value = 3838.7646
$ws.Range("I34").Value = 1792.8182
$ws.Range("J34").Value = 5767.8
$ws.Range("K34").Value = 1792.8182
$ws.Range("L34").Value = 5767.8
$ws.Range("M34").Value = -1590.8182
$ws.Range("N34").Value = -6171.8

$ws.Range("H51").Value = 79999.5
$ws.Range("J51").Value = 79999.5
$ws.Range("L51").Value = 79999.5
$ws.Range("N51").Value = -81471.5

$ws.Range("H58").Value = 31251960
$ws.Range("I58").Value = 45455816
$ws.Range("J58").Value = 3474
$ws.Range("K58").Value = 45455816
$ws.Range("L58").Value = 3474
$ws.Range("M58").Value = -45455613
$ws.Range("N58").Value = -3880

$ws.Range("H61").Value = 79999.5
$ws.Range("J61").Value = 79999.5
$ws.Range("L61").Value = 79999.5
$ws.Range("N61").Value = -80695.5

$ws.Range("H134").Value = 40543116
$ws.Range("I134").Value = 43479884
$ws.Range("J134").Value = 35718420
$ws.Range("K134").Value = 130439652
$ws.Range("L134").Value = 107155260
$ws.Range("M134").Value = -130437117
$ws.Range("N134").Value = -107160330

$ws.Range("H136").Value = 31251960
$ws.Range("I136").Value = 45455816
$ws.Range("J136").Value = 3474
$ws.Range("K136").Value = 136367448
$ws.Range("L136").Value = 10422
$ws.Range("M136").Value = -136364898
$ws.Range("N136").Value = -15522

$ws.Range("H138").Value = 49095
$ws.Range("J138").Value = 49095
$ws.Range("L138").Value = 49095
$ws.Range("N138").Value = -59375

$ws.Range("H141").Value = 241610.83
$ws.Range("J141").Value = 241610.83
$ws.Range("L141").Value = 241610.83
$ws.Range("N141").Value = -251970.83

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10102966
$ws.Range("J131").Value = 11496357
$ws.Range("L131").Value = 34489071
$ws.Range("N131").Value = -34499151

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = ""

$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -6840

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = ""

$ws.Range("H132").Value = 2814.919
$ws.Range("I132").Value = 2576.7083
$ws.Range("J132").Value = 3254.6924
$ws.Range("K132").Value = 7730.124899999999
$ws.Range("L132").Value = 9764.0772
$ws.Range("M132").Value = -5200.124899999999
$ws.Range("N132").Value = -14824.0772

$ws.Range("H133").Value = 52277.668
$ws.Range("J133").Value = 52277.668
$ws.Range("L133").Value = 52277.668
$ws.Range("N133").Value = -62397.668

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 3356327
$ws.Range("J81").Value = 3356327
$ws.Range("L81").Value = 3356327
$ws.Range("N81").Value = -3358323

$ws.Range("H84").Value = 3356327
$ws.Range("J84").Value = 3356327
$ws.Range("L84").Value = 10068981
$ws.Range("N84").Value = -10078965

$ws.Range("H141").Value = 61200
$ws.Range("J141").Value = 61200
$ws.Range("L141").Value = 61200
$ws.Range("N141").Value = -71560
